$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 125
$ws.Range("F2").Value = 89
$ws.Range("H2").Value = 95

$ws.Range("E3").Value = 51
$ws.Range("F3").Value = 37
$ws.Range("H3").Value = 38

$ws.Range("E5").Value = 178

$ws.Range("E10").Value = 809
$ws.Range("F10").Value = 476
$ws.Range("H10").Value = 571

$ws.Range("E11").Value = 539
$ws.Range("F11").Value = 325
$ws.Range("H11").Value = 390

$ws.Range("E12").Value = 834
$ws.Range("F12").Value = 525
$ws.Range("H12").Value = 611

$ws.Range("E13").Value = 184
$ws.Range("F13").Value = 108
$ws.Range("H13").Value = 142

$ws.Range("E14").Value = 164
$ws.Range("F14").Value = 94
$ws.Range("H14").Value = 128

$ws.Range("E15").Value = 224
$ws.Range("F15").Value = 109
$ws.Range("H15").Value = 160

$ws.Range("E16").Value = 252
$ws.Range("F16").Value = 150
$ws.Range("H16").Value = 198

$ws.Range("E17").Value = 140

$ws.Range("E18").Value = 69

$ws.Range("E20").Value = 114
$ws.Range("F20").Value = 54
$ws.Range("H20").Value = 91

$ws.Range("E21").Value = 162
$ws.Range("F21").Value = 101
$ws.Range("H21").Value = 132

$ws.Range("E22").Value = 217
$ws.Range("F22").Value = 126
$ws.Range("H22").Value = 168

$ws.Range("E23").Value = 256
$ws.Range("F23").Value = 132
$ws.Range("H23").Value = 184

$ws.Range("E24").Value = 318
$ws.Range("F24").Value = 184
$ws.Range("H24").Value = 214

$ws.Range("E25").Value = 381
$ws.Range("F25").Value = 218
$ws.Range("H25").Value = 278

$ws.Range("E26").Value = 247
$ws.Range("F26").Value = 149
$ws.Range("H26").Value = 174

$ws.Range("E27").Value = 440
$ws.Range("F27").Value = 255
$ws.Range("H27").Value = 337

$ws.Range("E28").Value = 255
$ws.Range("F28").Value = 128
$ws.Range("H28").Value = 180

$ws.Range("E29").Value = 207
$ws.Range("F29").Value = 130
$ws.Range("H29").Value = 171

$ws.Range("E30").Value = 290
$ws.Range("F30").Value = 183
$ws.Range("H30").Value = 236

$ws.Range("E31").Value = 89

$ws.Range("E32").Value = 239
$ws.Range("F32").Value = 157
$ws.Range("H32").Value = 195

$ws.Range("E33").Value = 373
$ws.Range("F33").Value = 201
$ws.Range("H33").Value = 292

$ws.Range("E34").Value = 283
$ws.Range("F34").Value = 201
$ws.Range("H34").Value = 239

$ws.Range("E35").Value = 197
$ws.Range("F35").Value = 141
$ws.Range("H35").Value = 168

$ws.Range("E37").Value = 215
$ws.Range("F37").Value = 122
$ws.Range("H37").Value = 158

$ws.Range("E38").Value = 114

$ws.Range("E39").Value = 219
$ws.Range("F39").Value = 114
$ws.Range("H39").Value = 165

$ws.Range("E40").Value = 339
$ws.Range("F40").Value = 187
$ws.Range("H40").Value = 267

$ws.Range("E41").Value = 491
$ws.Range("F41").Value = 264
$ws.Range("H41").Value = 356

$ws.Range("E42").Value = 519
$ws.Range("F42").Value = 321
$ws.Range("H42").Value = 382

$ws.Range("E43").Value = 160

$ws.Range("E44").Value = 435
$ws.Range("F44").Value = 248
$ws.Range("H44").Value = 316

$ws.Range("E45").Value = 205
$ws.Range("F45").Value = 124
$ws.Range("H45").Value = 163

$ws.Range("E46").Value = 432
$ws.Range("F46").Value = 259
$ws.Range("H46").Value = 323

$ws.Range("E47").Value = 621
$ws.Range("F47").Value = 365
$ws.Range("H47").Value = 457

$ws.Range("E48").Value = 304
$ws.Range("F48").Value = 156
$ws.Range("H48").Value = 200

$ws.Range("E49").Value = 373
$ws.Range("F49").Value = 193
$ws.Range("H49").Value = 280

$ws.Range("E50").Value = 309
$ws.Range("F50").Value = 178
$ws.Range("H50").Value = 251

$ws.Range("E51").Value = 277
$ws.Range("F51").Value = 145
$ws.Range("H51").Value = 219

$ws.Range("E52").Value = 36
$ws.Range("F52").Value = 17
$ws.Range("H52").Value = 25
